# Added get_form_type in the import_utils
#
# Adds a new "Form Tag" column (column V) to the CapitalCommitment sheet,
# with header "Form Tag" and a "Default" value filled in for every existing
# data row (rows 2-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added column.
$ws.Range("V1").Value = "Form Tag"

# Populate the new column for every existing data row with "Default".
$ws.Range("V2:V8").Value = "Default"

# Leave the new column selected, matching the state the workbook was left in
# after the interactive edit.
$ws.Range("V2:V8").Select()
